# Update the PROTECTION Okanogan steelhead habitat-quality sheet:
#  - rework the header row (renamed/reordered score columns, new summary
#    columns HQ_Sum/HQ_Pct/HQ_Score_Restoration/HQ_Score_Protection)
#  - update per-reach values, drop the now-unused reach-flag columns
#    (Spring.Chinook.Reach / Bull.Trout.Reach), and replace the stray
#    "Inf" temperature placeholders with real numeric data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Header row (row 1) --------------------------------------------------
$ws.Range("G1").Value = "BankStability_score"
$ws.Range("I1").Value = "Stability_Mean"
$ws.Range("J1").Value = "CoarseSubstrate_score"
$ws.Range("K1").Value = "Cover-Wood_score"
$ws.Range("L1").Value = "Flow-SummerBaseFlow_score"
$ws.Range("M1").Value = "Off-Channel-Floodplain_score"
$ws.Range("N1").Value = "Off-Channel-Side-Channels_score"
$ws.Range("O1").Value = "PoolQuantity&Quality_score"
$ws.Range("P1").Value = "Riparian-CanopyCover_score"
$ws.Range("Q1").Value = "Riparian-Disturbance_score"
$ws.Range("R1").Value = "Riparian_Mean"
$ws.Range("S1").Value = "Temperature-Rearing_score"
$ws.Range("T1").Value = "HQ_Sum"
$ws.Range("U1").Value = "HQ_Pct"
$ws.Range("V1").Value = "HQ_Score_Restoration"
$ws.Range("W1").Value = "HQ_Score_Protection"

# ---- Data rows 2-5 --------------------------------------------------------
# Columns D (Spring.Chinook.Reach) and F (Bull.Trout.Reach) are no longer
# populated for these reaches, and the old J/M/N/O/P columns are superseded
# by the reworked layout above, so clear anything not explicitly set below.
$dataRows = 2,3,4,5
foreach ($r in $dataRows) {
    foreach ($col in "D","F","J","M","N","O","P","T") {
        $ws.Range("$col$r").Value = $null
    }
}

# Row 2 - Salmon 16-11
$ws.Range("C2").Value = "Salmon Creek-Green Lake"
$ws.Range("G2").Value = 5
$ws.Range("H2").Value = 5
$ws.Range("I2").Value = 5
$ws.Range("K2").Value = 5
$ws.Range("L2").Value = 1
$ws.Range("Q2").Value = 3
$ws.Range("R2").Value = 3
$ws.Range("S2").Value = 5
$ws.Range("U2").Value = 0.872302720070069
$ws.Range("V2").Value = 3
$ws.Range("W2").Value = 3

# Row 3 - Salmon 16-6
$ws.Range("C3").Value = "Salmon Creek-Green Lake"
$ws.Range("G3").Value = 5
$ws.Range("H3").Value = 5
$ws.Range("I3").Value = 5
$ws.Range("K3").Value = 5
$ws.Range("L3").Value = 1
$ws.Range("Q3").Value = 1
$ws.Range("R3").Value = 1
$ws.Range("S3").Value = 5
$ws.Range("U3").Value = 0.7830164640163568
$ws.Range("V3").Value = 5
$ws.Range("W3").Value = 3

# Row 4 - Salmon 16-9
$ws.Range("C4").Value = "Salmon Creek-Green Lake"
$ws.Range("G4").Value = 5
$ws.Range("H4").Value = 5
$ws.Range("I4").Value = 5
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("Q4").Value = 1
$ws.Range("R4").Value = 1
$ws.Range("S4").Value = 5
$ws.Range("U4").Value = 0.7892917935080431
$ws.Range("V4").Value = 5
$ws.Range("W4").Value = 3

# Row 5 - Tonasket 16-2
$ws.Range("C5").Value = "Tonasket Creek"
$ws.Range("G5").Value = 5
$ws.Range("H5").Value = 5
$ws.Range("I5").Value = 5
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 3
$ws.Range("Q5").Value = 3
$ws.Range("R5").Value = 3
$ws.Range("S5").Value = 3
$ws.Range("U5").Value = 1
$ws.Range("V5").Value = 1
$ws.Range("W5").Value = 5
